$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$r2a = New-Object 'object[,]' 1,8
$r2a[0,0] = "name"
$r2a[0,1] = "anchor score"
$r2a[0,2] = "type occurences"
$r2a[0,3] = "total occurences"
$r2a[0,4] = "+%"
$r2a[0,5] = "-%"
$r2a[0,6] = "both"
$r2a[0,7] = "normal"
$ws.Range("A2:H2").Value = $r2a
$r2b = New-Object 'object[,]' 1,8
$r2b[0,0] = "name"
$r2b[0,1] = "anchor score"
$r2b[0,2] = "type occurences"
$r2b[0,3] = "total occurences"
$r2b[0,4] = "+%"
$r2b[0,5] = "-%"
$r2b[0,6] = "both"
$r2b[0,7] = "normal"
$ws.Range("J2:Q2").Value = $r2b

$r3a = New-Object 'object[,]' 1,8
$r3a[0,0] = "poorly"
$r3a[0,1] = 1
$r3a[0,2] = 46
$r3a[0,3] = 46
$r3a[0,4] = 0
$r3a[0,5] = 1
$r3a[0,6] = $false
$r3a[0,7] = 0
$ws.Range("A3:H3").Value = $r3a
$r3b = New-Object 'object[,]' 1,8
$r3b[0,0] = "wonderful"
$r3b[0,1] = 0.8571428571428571
$r3b[0,2] = 48
$r3b[0,3] = 48
$r3b[0,4] = 1
$r3b[0,5] = 0
$r3b[0,6] = $false
$r3b[0,7] = 8
$ws.Range("J3:Q3").Value = $r3b

$r4a = New-Object 'object[,]' 1,8
$r4a[0,0] = "disappointing"
$r4a[0,1] = 0.7954545454545454
$r4a[0,2] = 35
$r4a[0,3] = 35
$r4a[0,4] = 0
$r4a[0,5] = 1
$r4a[0,6] = $false
$r4a[0,7] = 9
$ws.Range("A4:H4").Value = $r4a
$r4b = New-Object 'object[,]' 1,8
$r4b[0,0] = "awesome"
$r4b[0,1] = 0.8461538461538461
$r4b[0,2] = 55
$r4b[0,3] = 55
$r4b[0,4] = 1
$r4b[0,5] = 0
$r4b[0,6] = $false
$r4b[0,7] = 10
$ws.Range("J4:Q4").Value = $r4b

$r5a = New-Object 'object[,]' 1,8
$r5a[0,0] = "however"
$r5a[0,1] = 0.78125
$r5a[0,2] = 50
$r5a[0,3] = 50
$r5a[0,4] = 0
$r5a[0,5] = 1
$r5a[0,6] = $false
$r5a[0,7] = 14
$ws.Range("A5:H5").Value = $r5a
$r5b = New-Object 'object[,]' 1,8
$r5b[0,0] = "favorite"
$r5b[0,1] = 0.7204301075268817
$r5b[0,2] = 67
$r5b[0,3] = 67
$r5b[0,4] = 1
$r5b[0,5] = 0
$r5b[0,6] = $false
$r5b[0,7] = 26
$ws.Range("J5:Q5").Value = $r5b

$r6a = New-Object 'object[,]' 1,8
$r6a[0,0] = "broke"
$r6a[0,1] = 0.7766990291262136
$r6a[0,2] = 160
$r6a[0,3] = 160
$r6a[0,4] = 0
$r6a[0,5] = 1
$r6a[0,6] = $false
$r6a[0,7] = 46
$ws.Range("A6:H6").Value = $r6a
$r6b = New-Object 'object[,]' 1,8
$r6b[0,0] = "classic"
$r6b[0,1] = 0.5660377358490566
$r6b[0,2] = 30
$r6b[0,3] = 30
$r6b[0,4] = 1
$r6b[0,5] = 0
$r6b[0,6] = $false
$r6b[0,7] = 23
$ws.Range("J6:Q6").Value = $r6b

$r7a = New-Object 'object[,]' 1,8
$r7a[0,0] = "disappointed"
$r7a[0,1] = 0.7580645161290323
$r7a[0,2] = 141
$r7a[0,3] = 141
$r7a[0,4] = 0
$r7a[0,5] = 1
$r7a[0,6] = $false
$r7a[0,7] = 45
$ws.Range("A7:H7").Value = $r7a
$r7b = New-Object 'object[,]' 1,8
$r7b[0,0] = "excellent"
$r7b[0,1] = 0.515625
$r7b[0,2] = 33
$r7b[0,3] = 33
$r7b[0,4] = 1
$r7b[0,5] = 0
$r7b[0,6] = $false
$r7b[0,7] = 31
$ws.Range("J7:Q7").Value = $r7b

$r8a = New-Object 'object[,]' 1,8
$r8a[0,0] = "poor"
$r8a[0,1] = 0.6901408450704225
$r8a[0,2] = 49
$r8a[0,3] = 49
$r8a[0,4] = 0
$r8a[0,5] = 1
$r8a[0,6] = $false
$r8a[0,7] = 22
$ws.Range("A8:H8").Value = $r8a
$r8b = New-Object 'object[,]' 1,8
$r8b[0,0] = "thank"
$r8b[0,1] = 0.4347826086956522
$r8b[0,2] = 30
$r8b[0,3] = 30
$r8b[0,4] = 1
$r8b[0,5] = 0
$r8b[0,6] = $false
$r8b[0,7] = 39
$ws.Range("J8:Q8").Value = $r8b

$r9a = New-Object 'object[,]' 1,8
$r9a[0,0] = "junk"
$r9a[0,1] = 0.6545454545454545
$r9a[0,2] = 36
$r9a[0,3] = 36
$r9a[0,4] = 0
$r9a[0,5] = 1
$r9a[0,6] = $false
$r9a[0,7] = 19
$ws.Range("A9:H9").Value = $r9a
$r9b = New-Object 'object[,]' 1,8
$r9b[0,0] = "great"
$r9b[0,1] = 0.3434426229508197
$r9b[0,2] = 419
$r9b[0,3] = 419
$r9b[0,4] = 1
$r9b[0,5] = 0
$r9b[0,6] = $false
$r9b[0,7] = 801
$ws.Range("J9:Q9").Value = $r9b

$r10a = New-Object 'object[,]' 1,8
$r10a[0,0] = "instead"
$r10a[0,1] = 0.6458333333333334
$r10a[0,2] = 31
$r10a[0,3] = 31
$r10a[0,4] = 0
$r10a[0,5] = 1
$r10a[0,6] = $false
$r10a[0,7] = 17
$ws.Range("A10:H10").Value = $r10a
$r10b = New-Object 'object[,]' 1,8
$r10b[0,0] = "love"
$r10b[0,1] = 0.2959770114942529
$r10b[0,2] = 206
$r10b[0,3] = 207
$r10b[0,4] = 1
$r10b[0,5] = 0
$r10b[0,6] = $true
$r10b[0,7] = 490
$ws.Range("J10:Q10").Value = $r10b

$r11a = New-Object 'object[,]' 1,8
$r11a[0,0] = "waste"
$r11a[0,1] = 0.6351351351351351
$r11a[0,2] = 94
$r11a[0,3] = 94
$r11a[0,4] = 0
$r11a[0,5] = 1
$r11a[0,6] = $false
$r11a[0,7] = 54
$ws.Range("A11:H11").Value = $r11a
$r11b = New-Object 'object[,]' 1,8
$r11b[0,0] = "loves"
$r11b[0,1] = 0.2634854771784232
$r11b[0,2] = 127
$r11b[0,3] = 127
$r11b[0,4] = 1
$r11b[0,5] = 0
$r11b[0,6] = $false
$r11b[0,7] = 355
$ws.Range("J11:Q11").Value = $r11b

$r12a = New-Object 'object[,]' 1,8
$r12a[0,0] = "smaller"
$r12a[0,1] = 0.5966386554621849
$r12a[0,2] = 71
$r12a[0,3] = 71
$r12a[0,4] = 0
$r12a[0,5] = 1
$r12a[0,6] = $false
$r12a[0,7] = 48
$ws.Range("A12:H12").Value = $r12a
$r12b = New-Object 'object[,]' 1,8
$r12b[0,0] = "loved"
$r12b[0,1] = 0.1834862385321101
$r12b[0,2] = 60
$r12b[0,3] = 60
$r12b[0,4] = 1
$r12b[0,5] = 0
$r12b[0,6] = $false
$r12b[0,7] = 267
$ws.Range("J12:Q12").Value = $r12b

$r13a = New-Object 'object[,]' 1,8
$r13a[0,0] = "paint"
$r13a[0,1] = 0.5079365079365079
$r13a[0,2] = 32
$r13a[0,3] = 32
$r13a[0,4] = 0
$r13a[0,5] = 1
$r13a[0,6] = $false
$r13a[0,7] = 31
$ws.Range("A13:H13").Value = $r13a
$r13b = New-Object 'object[,]' 1,8
$r13b[0,0] = "perfect"
$r13b[0,1] = 0.1746987951807229
$r13b[0,2] = 29
$r13b[0,3] = 29
$r13b[0,4] = 1
$r13b[0,5] = 0
$r13b[0,6] = $false
$r13b[0,7] = 137
$ws.Range("J13:Q13").Value = $r13b

$r14a = New-Object 'object[,]' 1,8
$r14a[0,0] = "small"
$r14a[0,1] = 0.4985507246376812
$r14a[0,2] = 172
$r14a[0,3] = 172
$r14a[0,4] = 0
$r14a[0,5] = 1
$r14a[0,6] = $false
$r14a[0,7] = 173
$ws.Range("A14:H14").Value = $r14a
$r14b = New-Object 'object[,]' 1,8
$r14b[0,0] = "fun"
$r14b[0,1] = 0.08326029798422437
$r14b[0,2] = 95
$r14b[0,3] = 95
$r14b[0,4] = 1
$r14b[0,5] = 0
$r14b[0,6] = $false
$r14b[0,7] = 1046
$ws.Range("J14:Q14").Value = $r14b

$r15a = New-Object 'object[,]' 1,8
$r15a[0,0] = "plastic"
$r15a[0,1] = 0.4409448818897638
$r15a[0,2] = 56
$r15a[0,3] = 56
$r15a[0,4] = 0
$r15a[0,5] = 1
$r15a[0,6] = $false
$r15a[0,7] = 71
$ws.Range("A15:H15").Value = $r15a
$r15b = New-Object 'object[,]' 1,8
$r15b[0,0] = "game"
$r15b[0,1] = 0.03893575600259572
$r15b[0,2] = 60
$r15b[0,3] = 60
$r15b[0,4] = 1
$r15b[0,5] = 0
$r15b[0,6] = $false
$r15b[0,7] = 1481
$ws.Range("J15:Q15").Value = $r15b

$r16a = New-Object 'object[,]' 1,8
$r16a[0,0] = "apart"
$r16a[0,1] = 0.4315789473684211
$r16a[0,2] = 41
$r16a[0,3] = 41
$r16a[0,4] = 0
$r16a[0,5] = 1
$r16a[0,6] = $false
$r16a[0,7] = 54
$ws.Range("A16:H16").Value = $r16a

$r17a = New-Object 'object[,]' 1,8
$r17a[0,0] = "broken"
$r17a[0,1] = 0.4096385542168675
$r17a[0,2] = 34
$r17a[0,3] = 34
$r17a[0,4] = 0
$r17a[0,5] = 1
$r17a[0,6] = $false
$r17a[0,7] = 49
$ws.Range("A17:H17").Value = $r17a

$r18a = New-Object 'object[,]' 1,8
$r18a[0,0] = "ok"
$r18a[0,1] = 0.3203125
$r18a[0,2] = 41
$r18a[0,3] = 41
$r18a[0,4] = 0
$r18a[0,5] = 1
$r18a[0,6] = $false
$r18a[0,7] = 87
$ws.Range("A18:H18").Value = $r18a

$r19a = New-Object 'object[,]' 1,8
$r19a[0,0] = "though"
$r19a[0,1] = 0.2735042735042735
$r19a[0,2] = 32
$r19a[0,3] = 32
$r19a[0,4] = 0
$r19a[0,5] = 1
$r19a[0,6] = $false
$r19a[0,7] = 85
$ws.Range("A19:H19").Value = $r19a

$r20a = New-Object 'object[,]' 1,8
$r20a[0,0] = "thought"
$r20a[0,1] = 0.2623762376237624
$r20a[0,2] = 53
$r20a[0,3] = 53
$r20a[0,4] = 0
$r20a[0,5] = 1
$r20a[0,6] = $false
$r20a[0,7] = 149
$ws.Range("A20:H20").Value = $r20a

$r21a = New-Object 'object[,]' 1,8
$r21a[0,0] = "cheap"
$r21a[0,1] = 0.2559241706161137
$r21a[0,2] = 54
$r21a[0,3] = 54
$r21a[0,4] = 0
$r21a[0,5] = 1
$r21a[0,6] = $false
$r21a[0,7] = 157
$ws.Range("A21:H21").Value = $r21a

$r22a = New-Object 'object[,]' 1,8
$r22a[0,0] = "size"
$r22a[0,1] = 0.2422680412371134
$r22a[0,2] = 47
$r22a[0,3] = 47
$r22a[0,4] = 0
$r22a[0,5] = 1
$r22a[0,6] = $false
$r22a[0,7] = 147
$ws.Range("A22:H22").Value = $r22a

$r23a = New-Object 'object[,]' 1,8
$r23a[0,0] = "work"
$r23a[0,1] = 0.1962025316455696
$r23a[0,2] = 62
$r23a[0,3] = 62
$r23a[0,4] = 0
$r23a[0,5] = 1
$r23a[0,6] = $false
$r23a[0,7] = 254
$ws.Range("A23:H23").Value = $r23a

$r24a = New-Object 'object[,]' 1,8
$r24a[0,0] = "item"
$r24a[0,1] = 0.1884057971014493
$r24a[0,2] = 52
$r24a[0,3] = 52
$r24a[0,4] = 0
$r24a[0,5] = 1
$r24a[0,6] = $false
$r24a[0,7] = 224
$ws.Range("A24:H24").Value = $r24a

$r25a = New-Object 'object[,]' 1,8
$r25a[0,0] = "money"
$r25a[0,1] = 0.180379746835443
$r25a[0,2] = 57
$r25a[0,3] = 57
$r25a[0,4] = 0
$r25a[0,5] = 1
$r25a[0,6] = $false
$r25a[0,7] = 259
$ws.Range("A25:H25").Value = $r25a

$r26a = New-Object 'object[,]' 1,8
$r26a[0,0] = "would"
$r26a[0,1] = 0.172106824925816
$r26a[0,2] = 116
$r26a[0,3] = 116
$r26a[0,4] = 0
$r26a[0,5] = 1
$r26a[0,6] = $false
$r26a[0,7] = 558
$ws.Range("A26:H26").Value = $r26a

$r27a = New-Object 'object[,]' 1,8
$r27a[0,0] = "hard"
$r27a[0,1] = 0.16
$r27a[0,2] = 32
$r27a[0,3] = 32
$r27a[0,4] = 0
$r27a[0,5] = 1
$r27a[0,6] = $false
$r27a[0,7] = 168
$ws.Range("A27:H27").Value = $r27a

$r28a = New-Object 'object[,]' 1,8
$r28a[0,0] = "better"
$r28a[0,1] = 0.1448598130841121
$r28a[0,2] = 31
$r28a[0,3] = 31
$r28a[0,4] = 0
$r28a[0,5] = 1
$r28a[0,6] = $false
$r28a[0,7] = 183
$ws.Range("A28:H28").Value = $r28a

$r29a = New-Object 'object[,]' 1,8
$r29a[0,0] = "product"
$r29a[0,1] = 0.1343612334801762
$r29a[0,2] = 61
$r29a[0,3] = 61
$r29a[0,4] = 0
$r29a[0,5] = 1
$r29a[0,6] = $false
$r29a[0,7] = 393
$ws.Range("A29:H29").Value = $r29a

$r30a = New-Object 'object[,]' 1,8
$r30a[0,0] = "price"
$r30a[0,1] = 0.1206896551724138
$r30a[0,2] = 42
$r30a[0,3] = 42
$r30a[0,4] = 0
$r30a[0,5] = 1
$r30a[0,6] = $false
$r30a[0,7] = 306
$ws.Range("A30:H30").Value = $r30a

$r31a = New-Object 'object[,]' 1,8
$r31a[0,0] = "2"
$r31a[0,1] = 0.1086142322097378
$r31a[0,2] = 29
$r31a[0,3] = 29
$r31a[0,4] = 0
$r31a[0,5] = 1
$r31a[0,6] = $false
$r31a[0,7] = 238
$ws.Range("A31:H31").Value = $r31a

$r32a = New-Object 'object[,]' 1,8
$r32a[0,0] = "use"
$r32a[0,1] = 0.09863013698630137
$r32a[0,2] = 36
$r32a[0,3] = 36
$r32a[0,4] = 0
$r32a[0,5] = 1
$r32a[0,6] = $false
$r32a[0,7] = 329
$ws.Range("A32:H32").Value = $r32a

$r33a = New-Object 'object[,]' 1,8
$r33a[0,0] = "like"
$r33a[0,1] = 0.066006600660066
$r33a[0,2] = 40
$r33a[0,3] = 42
$r33a[0,4] = 0.05
$r33a[0,5] = 0.95
$r33a[0,6] = $true
$r33a[0,7] = 566
$ws.Range("A33:H33").Value = $r33a

$r34a = New-Object 'object[,]' 1,8
$r34a[0,0] = "little"
$r34a[0,1] = 0.0645879732739421
$r34a[0,2] = 29
$r34a[0,3] = 29
$r34a[0,4] = 0
$r34a[0,5] = 1
$r34a[0,6] = $false
$r34a[0,7] = 420
$ws.Range("A34:H34").Value = $r34a

# Cells whose text looks like a number need to be forced back to text
# (Excel auto-converts numeric-looking strings on assignment); re-enter
# each as a formula producing the text, then paste-values-only onto itself
# so the stored cell is a literal string again, matching the source file.
$ws.Range("A31").Formula = '="2"'
$ws.Range("A31").Copy()
$ws.Range("A31").PasteSpecial(-4163)
$excel.CutCopyMode = $false
